$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")
$ws.Activate()

# The "valueType" column (D) for the data rows (8-101) is being changed from
# "integer" to "decimal". Row 8's cell (D8) through row 15 already carry the
# formatting used for "decimal" values elsewhere on the sheet (e.g. the BMI
# row, D7), so copy that formatting onto the remaining cells (D16:D101) that
# still used the older "integer" style before overwriting all of the values.
$formatSource = $ws.Cells.Item(8, 4)
$formatSource.Copy()
$ws.Range("D16:D101").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 8; $r -le 101; $r++) {
    $ws.Cells.Item($r, 4).Value = "decimal"
}

# Reflect the selection the user left on the sheet after performing the bulk
# edit on the valueType column.
[void]$ws.Range("D8:D101").Select()
